# --- Cell values (DB schema note) ---
# order chosen to match the original authoring sequence (shared-string order)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "Salt"
$ws.Range("N4").Value = "Hash-Algo"
$ws.Range("A2").Value = "DB"
$ws.Range("T2").Value = "Function"
$ws.Range("A6").Value = "password"
$ws.Range("A7").Value = "authentication"
$ws.Range("B6").Value = "varchar"
$ws.Range("A3").Value = "id"
$ws.Range("B3").Value = "int"
$ws.Range("B7").Value = "tinyInt/bool"
$ws.Range("B2").Value = "user"
$ws.Range("A8").Value = "SecureWord"
$ws.Range("B8").Value = "varchar"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 11.7109375
$ws.Columns.Item(14).ColumnWidth = 10

# --- Selection ---
$ws.Range("A7").Select()
